# Apply the "rallies" sheet update:
#  - D85 changes from 22 to 23 (rally_no for the last existing rally)
#  - A new row 86 is appended, duplicating row 85's data but for the
#    new rally (rally_id=85, rally_no=23, score_home=23)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rallies")

# Update existing cell D85: rally_no 22 -> 23
$ws.Range("D85").Value = 23

# Append new row 86 with the new rally's data
$ws.Range("A86").Value = 85
$ws.Range("B86").Value = 1
$ws.Range("C86").Value = 3
$ws.Range("D86").Value = 23
$ws.Range("E86").Value = "NOS"
$ws.Range("F86").Value = ""
$ws.Range("G86").Value = 5
$ws.Range("H86").Value = "LOB"
$ws.Range("I86").Value = "PONTO"
$ws.Range("J86").Value = "NOS"
$ws.Range("K86").Value = 23
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = "1 5 lob"
$ws.Range("N86").Value = "FRENTE"
$ws.Range("O86").Value = "FRENTE"
$ws.Range("P86").Value = "FRENTE"
